$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 307, shifting existing rows 307:400 down to 308:401.
$ws.Rows.Item(307).Insert()

# Populate the newly inserted row 307 with a new weekly data point
# (same dimensions/metadata as its neighbours, new date + volume).
$ws.Cells.Item(307, 1).Value = 10
$ws.Cells.Item(307, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(307, 3).Value = "La Araucanía"
$ws.Cells.Item(307, 4).Value = 45120
$ws.Cells.Item(307, 5).Value = 9
$ws.Cells.Item(307, 6).Value = 100112039
$ws.Cells.Item(307, 7).Value = "Ciboulette"
$ws.Cells.Item(307, 8).Value = "Sin especificar"
$ws.Cells.Item(307, 9).Value = "Primera"
$ws.Cells.Item(307, 10).Value = 60
$ws.Cells.Item(307, 11).Value = 5000
$ws.Cells.Item(307, 12).Value = 5000
$ws.Cells.Item(307, 13).Value = 5000
$ws.Cells.Item(307, 14).Value = "$/docena de atados"
$ws.Cells.Item(307, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(307, 16).Value = 1667
$ws.Cells.Item(307, 17).Value = 3
$ws.Cells.Item(307, 18).Value = "Hortaliza"
